# Updated cryptos list with latest prices, volumes, and rankings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.324.65"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.427.25"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'581.27"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'177.93"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.418.93"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.197"
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.586"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'48.64"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000282"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'690.17"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.979.10"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'8.65"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.366.00"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.422.13"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.121"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'17.81"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'11.37"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.904"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'17.06"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'5.34"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'101.12"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'3.89"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'9.67"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'33.70"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.75"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.92"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "'3.77"
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'564.63"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.05"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'58.08"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.631.47"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'34.95"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0735"
$ws.Range("E41").Value = "  +7.32%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.27"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0423"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.334"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'1.49"
$ws.Range("E46").Value = "  +9.99%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.64"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'130.84"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "'2.70"
$ws.Range("E51").Value = "  +2.75%  "
